# [ADD] the repository of Predefined template use github
#
# Applies the 2018-5-3-devops.xlsx edit:
#  - repo_url / name cell text updated to real GitHub URLs for the
#    predefined app templates (MicroService, MicroServiceFront, JavaLib)
#  - those repo_url cells turned into real hyperlinks
#  - row 10 grows a touch taller to fit the longer URL
#  - the devops_app_template sheet's view scrolls/selects a different cell
#  - column widths widen a bit on all three sheets

$wb = $excel.ActiveWorkbook

$readme = $wb.Worksheets.Item("README")
$app    = $wb.Worksheets.Item("devops_app_template")
$user   = $wb.Worksheets.Item("devops_user")

# ---------------------------------------------------------------------
# 1. Cell text content -> new GitHub urls / renamed template
# ---------------------------------------------------------------------
$app.Range("G8").Value2 = "https://github.com/choerodon/choerodon-front-template.git"

$app.Range("E9").Value2 = "MicroServiceFront"
$app.Range("F9").Value2 = "MicroServiceFront"
$app.Range("H9").Value2 = "MicroServiceFront"
$app.Range("G9").Value2 = "https://github.com/choerodon/choerodon-microservice-template.git"

$app.Range("G10").Value2 = "https://github.com/choerodon/choerodon-javalib-template.git"

# ---------------------------------------------------------------------
# 2. Turn the repo_url cells into real hyperlinks
# ---------------------------------------------------------------------
$app.Hyperlinks.Add($app.Range("G8"), "https://github.com/choerodon/choerodon-front-template.git", "", "", "https://github.com/choerodon/choerodon-front-template.git")
$app.Hyperlinks.Add($app.Range("G9"), "https://github.com/choerodon/choerodon-microservice-template.git", "", "", "https://github.com/choerodon/choerodon-microservice-template.git")
$app.Hyperlinks.Add($app.Range("G10"), "https://github.com/choerodon/choerodon-javalib-template.git", "", "", "https://github.com/choerodon/choerodon-javalib-template.git")

# ---------------------------------------------------------------------
# 3. Row 10 grows slightly taller (wrapped long url)
# ---------------------------------------------------------------------
$app.Rows.Item(10).RowHeight = 15.7

# ---------------------------------------------------------------------
# 4. devops_app_template sheet view: scroll + selection move
# ---------------------------------------------------------------------
$app.Activate()
$excel.ActiveWindow.TopLeftCell = $app.Range("G1")
$app.Range("H9").Select()

# ---------------------------------------------------------------------
# 5. Column widths widen (~13-14%) on all three sheets
# ---------------------------------------------------------------------
$readme.Range($readme.Cells.Item(1, 1), $readme.Cells.Item(1, 1)).EntireColumn.ColumnWidth = 17.857142857142858
$readme.Range($readme.Cells.Item(1, 2), $readme.Cells.Item(1, 1011)).EntireColumn.ColumnWidth = 16.571428571428573
$readme.Range($readme.Cells.Item(1, 1012), $readme.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 17

$app.Range($app.Cells.Item(1, 1), $app.Cells.Item(1, 3)).EntireColumn.ColumnWidth = 16.571428571428573
$app.Range($app.Cells.Item(1, 4), $app.Cells.Item(1, 4)).EntireColumn.ColumnWidth = 32.714285714285715
$app.Range($app.Cells.Item(1, 5), $app.Cells.Item(1, 5)).EntireColumn.ColumnWidth = 54
$app.Range($app.Cells.Item(1, 6), $app.Cells.Item(1, 6)).EntireColumn.ColumnWidth = 37.142857142857146
$app.Range($app.Cells.Item(1, 7), $app.Cells.Item(1, 7)).EntireColumn.ColumnWidth = 232.28571428571428
$app.Range($app.Cells.Item(1, 8), $app.Cells.Item(1, 8)).EntireColumn.ColumnWidth = 30.428571428571427
$app.Range($app.Cells.Item(1, 9), $app.Cells.Item(1, 9)).EntireColumn.ColumnWidth = 39
$app.Range($app.Cells.Item(1, 10), $app.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 16.571428571428573

$user.Range($user.Cells.Item(1, 1), $user.Cells.Item(1, 3)).EntireColumn.ColumnWidth = 14.285714285714286
$user.Range($user.Cells.Item(1, 4), $user.Cells.Item(1, 4)).EntireColumn.ColumnWidth = 27.428571428571427
$user.Range($user.Cells.Item(1, 5), $user.Cells.Item(1, 5)).EntireColumn.ColumnWidth = 14.285714285714286
$user.Range($user.Cells.Item(1, 6), $user.Cells.Item(1, 6)).EntireColumn.ColumnWidth = 45
$user.Range($user.Cells.Item(1, 7), $user.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 14.285714285714286
